$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D cells contain numeric-looking text (prices). Force them to stay
# text (matching the source inlineStr cells) by pre-setting a text number
# format, then after writing the value, paste the format from an untouched
# text cell (D4) back on top so the stored style index is unchanged.
$dRows = @(2,3,5,6,9,10,11,13,14,15,17,18,19,20,21,22,23,25,26,28,32,33,35,36,37,40,46,47,49,51)
foreach ($r in $dRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = "75.479.65"
$ws.Range("E2").Value = "  +2.05%  "

$ws.Range("D3").Value = "2.823.27"
$ws.Range("E3").Value = "  +6.89%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "190.70"
$ws.Range("E5").Value = "  +2.65%  "

$ws.Range("D6").Value = "595.49"
$ws.Range("E6").Value = "  +2.32%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("E8").Value = "  +3.31%  "

$ws.Range("D9").Value = "0.192"
$ws.Range("E9").Value = "  +0.62%  "

$ws.Range("D10").Value = "2.823.53"
$ws.Range("E10").Value = "  +7.03%  "

$ws.Range("D11").Value = "0.381"
$ws.Range("E11").Value = "  +7.99%  "

$ws.Range("E12").Value = "  -1.95%  "

$ws.Range("D13").Value = "4.89"
$ws.Range("E13").Value = "  +4.77%  "

$ws.Range("D14").Value = "3.337.65"
$ws.Range("E14").Value = "  +5.89%  "

$ws.Range("D15").Value = "75.278.28"
$ws.Range("E15").Value = "  +1.77%  "

$ws.Range("E16").Value = "  +1.75%  "

$ws.Range("D17").Value = "26.98"
$ws.Range("E17").Value = "  +3.44%  "

$ws.Range("D18").Value = "2.820.36"
$ws.Range("E18").Value = "  +5.61%  "

$ws.Range("D19").Value = "9.02"
$ws.Range("E19").Value = "  -2.42%  "

$ws.Range("D20").Value = "12.32"
$ws.Range("E20").Value = "  +4.37%  "

$ws.Range("D21").Value = "378.81"
$ws.Range("E21").Value = "  +2.74%  "

$ws.Range("D22").Value = "2.31"
$ws.Range("E22").Value = "  +2.72%  "

$ws.Range("D23").Value = "4.10"
$ws.Range("E23").Value = "  +1.66%  "

$ws.Range("E24").Value = "  -0.07%  "

$ws.Range("D25").Value = "71.14"
$ws.Range("E25").Value = "  +2.15%  "

$ws.Range("B26").Value = "Aptos"
$ws.Range("C26").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D26").Value = "9.80"
$ws.Range("E26").Value = "  +6.22%  "

$ws.Range("E27").Value = "  +2.65%  "

$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.961.72"
$ws.Range("E28").Value = "  +5.90%  "

$ws.Range("E29").Value = "  +12.48%  "

$ws.Range("E30").Value = "  -0.27%  "

$ws.Range("E31").Value = "  +1.79%  "

$ws.Range("D32").Value = "515.42"
$ws.Range("E32").Value = "  +0.03%  "

$ws.Range("D33").Value = "7.70"
$ws.Range("E33").Value = "  +1.53%  "

$ws.Range("E34").Value = "  +4.47%  "

$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.12%  "

$ws.Range("D36").Value = "164.97"
$ws.Range("E36").Value = "  +1.34%  "

$ws.Range("D37").Value = "19.88"
$ws.Range("E37").Value = "  +4.23%  "

$ws.Range("E38").Value = "  +0.77%  "

$ws.Range("E39").Value = "  +0.45%  "

$ws.Range("D40").Value = "182.72"
$ws.Range("E40").Value = "  +11.48%  "

$ws.Range("E41").Value = "  -0.05%  "

$ws.Range("E42").Value = "  +5.65%  "

$ws.Range("E43").Value = "  +3.19%  "

$ws.Range("E44").Value = "  +1.72%  "

$ws.Range("E45").Value = "  +3.61%  "

$ws.Range("D46").Value = "39.97"
$ws.Range("E46").Value = "  +2.68%  "

$ws.Range("D47").Value = "0.0872"
$ws.Range("E47").Value = "  +3.53%  "

$ws.Range("E48").Value = "  +1.07%  "

$ws.Range("D49").Value = "0.571"
$ws.Range("E49").Value = "  +9.09%  "

$ws.Range("E50").Value = "  +4.15%  "

$ws.Range("D51").Value = "0.645"
$ws.Range("E51").Value = "  +9.94%  "

# Restore original (General/style-0) formatting on the column D cells we
# touched, copying it from D4 which is untouched by this update.
$ws.Range("D4").Copy() | Out-Null
foreach ($r in $dRows) {
    $ws.Range("D$r").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0
